# Auto-generated script to update "想去人数" (F column) values
# across all 4 worksheets per the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 245
$ws.Range("F4").Value = 516
$ws.Range("F5").Value = 2294
$ws.Range("F7").Value = 8170
$ws.Range("F8").Value = 116
$ws.Range("F10").Value = 1618
$ws.Range("F11").Value = 1329
$ws.Range("F13").Value = 4466
$ws.Range("F14").Value = 6176
$ws.Range("F15").Value = 785
$ws.Range("F17").Value = 1239
$ws.Range("F18").Value = 1268
$ws.Range("F19").Value = 483
$ws.Range("F20").Value = 6496
$ws.Range("F21").Value = 362
$ws.Range("F24").Value = 4393
$ws.Range("F25").Value = 317
$ws.Range("F26").Value = 721
$ws.Range("F27").Value = 2041
$ws.Range("F28").Value = 1192
$ws.Range("F29").Value = 353
$ws.Range("F30").Value = 1084
$ws.Range("F31").Value = 51
$ws.Range("F32").Value = 53
$ws.Range("F33").Value = 45
$ws.Range("F34").Value = 88
$ws.Range("F36").Value = 1198
$ws.Range("F37").Value = 1909
$ws.Range("F38").Value = 145
$ws.Range("F40").Value = 170
$ws.Range("F41").Value = 1217
$ws.Range("F42").Value = 565
$ws.Range("F43").Value = 71
$ws.Range("F44").Value = 1182
$ws.Range("F47").Value = 196
$ws.Range("F48").Value = 31
$ws.Range("F49").Value = 19

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 87
$ws.Range("F8").Value = 696
$ws.Range("F10").Value = 417
$ws.Range("F12").Value = 218
$ws.Range("F17").Value = 196
$ws.Range("F18").Value = 154
$ws.Range("F22").Value = 110
$ws.Range("F27").Value = 286
$ws.Range("F28").Value = 166

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 121
$ws.Range("F3").Value = 78
$ws.Range("F4").Value = 467
$ws.Range("F6").Value = 1593
$ws.Range("F8").Value = 3160
$ws.Range("F9").Value = 1075
$ws.Range("F10").Value = 1158
$ws.Range("F11").Value = 1497
$ws.Range("F12").Value = 1847
$ws.Range("F13").Value = 329
$ws.Range("F14").Value = 208

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 467
$ws.Range("F3").Value = 1593
$ws.Range("F5").Value = 516
$ws.Range("F7").Value = 3160
$ws.Range("F8").Value = 2294
$ws.Range("F9").Value = 1075
$ws.Range("F11").Value = 1618
$ws.Range("F12").Value = 1497
$ws.Range("F13").Value = 1329
$ws.Range("F14").Value = 696
$ws.Range("F16").Value = 1847
$ws.Range("F17").Value = 4466
$ws.Range("F18").Value = 329
$ws.Range("F20").Value = 417
$ws.Range("F21").Value = 785
$ws.Range("F23").Value = 1239
$ws.Range("F24").Value = 1268
$ws.Range("F25").Value = 483
$ws.Range("F26").Value = 6496
$ws.Range("F27").Value = 362
$ws.Range("F28").Value = 208
$ws.Range("F30").Value = 4393
$ws.Range("F31").Value = 317
$ws.Range("F32").Value = 2041
$ws.Range("F33").Value = 1192
$ws.Range("F34").Value = 353
$ws.Range("F35").Value = 1084
$ws.Range("F36").Value = 51
$ws.Range("F37").Value = 45
$ws.Range("F38").Value = 196
$ws.Range("F39").Value = 88
$ws.Range("F41").Value = 1909
$ws.Range("F42").Value = 145
$ws.Range("F44").Value = 1217
$ws.Range("F46").Value = 565
$ws.Range("F47").Value = 286
$ws.Range("F48").Value = 1182
$ws.Range("F49").Value = 196

$wb.Save()
Write-Output "Updated F column values across all sheets."
